$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 110 (new week's data for Ají, Primera/Segunda),
# pushing the existing rows 110-121 down to 112-123.
$ws.Range("A110:R111").EntireRow.Insert()

# New row 110: Primera
$ws.Range("A110").Value = 1
$ws.Range("B110").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C110").Value = "Arica y Parinacota"
$ws.Range("D110").Value = 44946
$ws.Range("E110").Value = 15
$ws.Range("F110").Value = 100112021
$ws.Range("G110").Value = "Ají"
$ws.Range("H110").Value = "Inferno"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 450
$ws.Range("K110").Value = 14000
$ws.Range("L110").Value = 15000
$ws.Range("M110").Value = 14444
$ws.Range("N110").Value = "$/caja 15 kilos"
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 963
$ws.Range("Q110").Value = 15
$ws.Range("R110").Value = "Hortaliza"

# New row 111: Segunda
$ws.Range("A111").Value = 1
$ws.Range("B111").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C111").Value = "Arica y Parinacota"
$ws.Range("D111").Value = 44946
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = 100112021
$ws.Range("G111").Value = "Ají"
$ws.Range("H111").Value = "Inferno"
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 300
$ws.Range("K111").Value = 11000
$ws.Range("L111").Value = 13000
$ws.Range("M111").Value = 12333
$ws.Range("N111").Value = "$/caja 15 kilos"
$ws.Range("O111").Value = "Región de Arica y Parinacota"
$ws.Range("P111").Value = 822
$ws.Range("Q111").Value = 15
$ws.Range("R111").Value = "Hortaliza"
